$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
Write-Host ($ws1 | Get-Member | Select-String -Pattern "Format")
